# Adding the final working code including the main python file with more tests
#
# Appends four new bulleted ("ListParagraph" / numId 1) paragraphs at the
# very end of the document, right before the final section break:
#   - "Create main file"                                                  (ilvl 0)
#   - "Pull data on all starships"                                        (ilvl 1)
#   - "Replace 'pilots' field with list of their respective ObjectIDs
#      from the characters collection"                                    (ilvl 1)
#   - "Insert the starships as documents into the starships collection
#      in MongoDB"                                                        (ilvl 1)

$d = $word.ActiveDocument
$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$paragraphsXml = @(
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Create main file</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Pull data on all </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>starships</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Replace ' + [char]0x2018 + 'pilots' + [char]0x2019 + ' field with list of their respective </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ObjectIDs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from the characters collection</w:t></w:r></w:p>'),
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Insert the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>starships</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as documents into the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>starships</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> collection in MongoDB</w:t></w:r></w:p>'
)

# Start a trailing scratch paragraph (inherits the last paragraph's
# formatting) that each new chunk of XML is inserted just ahead of; the
# scratch paragraph is always empty so InsertXML's "insert before the
# target range" behaviour never clobbers real content. The paragraph is
# re-fetched from the Paragraphs collection (rather than reusing the old
# Range reference) before every insertion so it tracks the live document.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null

foreach ($bodyXml in $paragraphsXml) {
    $scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $insertAt = $scratchPara.Range
    $insertAt.Collapse(1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        "<w:document $wordMlNs><w:body>$bodyXml</w:body></w:document>" +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $insertAt.InsertXML($pkg)
}

# Drop the still-empty scratch paragraph used as the perpetual insertion
# anchor above by deleting the paragraph mark that separates it from the
# real content before it (Range.Delete() on the scratch paragraph's own,
# pilcrow-only range is a no-op since it is the final paragraph mark of
# the document body).
$n = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($n - 1)
$d.Range($prevPara.Range.End - 1, $prevPara.Range.End).Delete()

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
